# Apply the "any comment here" annotation change described in the commit:
# both long protein-description annotations in the homolog rows (G7, G8)
# are replaced by a single shared placeholder comment string, and the
# active selection on Sheet1 moves to G18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G7").Value2 = "any comment here"
$ws.Range("G8").Value2 = "any comment here"

$ws.Range("G18").Select()
